$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve exact text formatting (e.g. trailing zeros) for price cells that
# would otherwise be auto-converted to numbers by Excel when set via .Value.
$textCells = 'D5','D6','D7','D10','D14','D19','D20','D23','D24','D26','D27','D28','D30','D32','D33','D37','D38','D40','D41','D42','D43','D44','D45','D47','D50'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '56.638.95'
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = '2.324.25'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '516.36'
$ws.Range('E5').Value = '  +1.98%  '
$ws.Range('D6').Value = '135.55'
$ws.Range('E6').Value = '  +5.54%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +1.50%  '
$ws.Range('D9').Value = '2.342.84'
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('D10').Value = '0.102'
$ws.Range('E10').Value = '  +3.64%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('E12').Value = '  +5.35%  '
$ws.Range('E13').Value = '  +0.37%  '
$ws.Range('D14').Value = '23.96'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '2.740.94'
$ws.Range('E15').Value = '  +2.24%  '
$ws.Range('D16').Value = '56.716.90'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '2.332.82'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').Value = '10.53'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').Value = '326.14'
$ws.Range('E20').Value = '  +3.69%  '
$ws.Range('E21').Value = '  +0.63%  '
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('D24').Value = '60.78'
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('E25').Value = '  +6.05%  '
$ws.Range('D26').Value = '0.995'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = '7.99'
$ws.Range('E27').Value = '  +5.84%  '
$ws.Range('D28').Value = '1.28'
$ws.Range('E28').Value = '  +11.25%  '
$ws.Range('D29').Value = '0.0₃0743'
$ws.Range('E29').Value = '  +5.13%  '
$ws.Range('D30').Value = '168.78'
$ws.Range('E30').Value = '  -1.37%  '
$ws.Range('E31').Value = '  +3.26%  '
$ws.Range('D32').Value = '6.20'
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('D33').Value = '18.47'
$ws.Range('E33').Value = '  +2.76%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('E36').Value = '  +1.99%  '
$ws.Range('D37').Value = '0.918'
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('D38').Value = '4.00'
$ws.Range('E38').Value = '  +2.65%  '
$ws.Range('E39').Value = '  +6.99%  '
$ws.Range('D40').Value = '38.32'
$ws.Range('E40').Value = '  +4.18%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '0.381'
$ws.Range('E41').Value = '  +1.84%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '142.08'
$ws.Range('E42').Value = '  +3.97%  '
$ws.Range('D43').Value = '3.60'
$ws.Range('E43').Value = '  +3.66%  '
$ws.Range('D44').Value = '5.21'
$ws.Range('E44').Value = '  +6.55%  '
$ws.Range('D45').Value = '278.28'
$ws.Range('E45').Value = '  +7.73%  '
$ws.Range('E46').Value = '  +1.81%  '
$ws.Range('D47').Value = '0.0507'
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('E48').Value = '  +2.65%  '
$ws.Range('E49').Value = '  +2.74%  '
$ws.Range('D50').Value = '17.87'
$ws.Range('E50').Value = '  +8.60%  '
$ws.Range('E51').Value = '  +1.61%  '

# Restore the default "Normal" style on the cells we had to temporarily
# mark as Text, so number formatting matches the original workbook.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
